# TODO update. Address the first block first!
#
# Rewrites the Sheet1 TODO list: removes the "9) זכרון" item, splits the
# "restart" wording out of item 7, adds a handful of new TODO rows (an
# rsa-based encryption note plus several AV/startup/macro/USB-drive checks
# sitting alongside the "OVERALL TEST" row), and shifts the whole list up
# so it now starts at row 2 and ends at row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe everything that is currently in use so stale rows/styles disappear.
$ws.Range("B3:C20").Clear()

# --- New content -----------------------------------------------------
$ws.Range("B2").Value2 = "rsa-based head-of-file encryption"
$ws.Range("B3").Value2 = "רפקטורינג למנגנון iterate"

$ws.Range("B4").Value2 = "15) חזרה בRESTART"
$ws.Range("B5").Value2 = "10) לעטוף את DEC ואת ENC לקובץ אחד"

$ws.Range("B8").Value2 = "7) הודעה ללקוח -  ביטול רישום ל-startup processes"

$ws.Range("B9").Value2 = "8) OVERALL TEST על מכונה וירטואלית"
$ws.Range("C9").Value2 = "כמה אנטיוירוסים"
$ws.Range("D9").Value2 = "בדיקת stratup"
$ws.Range("E9").Value2 = "בדיקת מאקרו"
$ws.Range("F9").Value2 = "6) BACKUP CORRUPTION"

$ws.Range("B10").Value2 = "לסגור ווינדוס דפנדר"
$ws.Range("B11").Value2 = "12) WORD  הנדסת אנוש"
$ws.Range("B12").Value2 = "מיפוי כוננים נוספים - כולל USB DRIVES"
$ws.Range("B13").Value2 = "13) סיומת קובץ 5 תווים אקריים/ קבצי PDF עם MAGIC"

$ws.Range("B15").Value2 = "גילוי קבצים מוצפנים בסיום הריצה"
$ws.Range("B16").Value2 = "14) להמנע מקבצים מוסתרים"
$ws.Range("B17").Value2 = "שינוי סדר איטרציה"
$ws.Range("B18").Value2 = "16) מיון קבצים נוסף לפי גודל"

# --- Styling (right aligned / RTL reading order style, same as before) --
$ws.Range("B4:B5").HorizontalAlignment = -4152
$ws.Range("B8:B13").HorizontalAlignment = -4152
$ws.Range("F9").HorizontalAlignment = -4152
$ws.Range("B15:B18").HorizontalAlignment = -4152
$ws.Range("C15").HorizontalAlignment = -4152  # C15 stays empty but keeps the style (mirrors old C14 placeholder)

# --- Column widths (best effort; engine rounds to whole pixels) --------
$ws.Columns(2).ColumnWidth = 42.666666666666664
$ws.Columns(3).ColumnWidth = 13.0
$ws.Columns(4).ColumnWidth = 11.666666666666666
$ws.Columns(5).ColumnWidth = 10.5

# --- Selection ----------------------------------------------------------
$ws.Range("B15").Select()
